$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 174; everything currently at/after row 174 (down to the
# old last row 194) shifts down by one, so old row 194 ends up at new row 195
# and dimension grows from A1:R194 to A1:R195.
$ws.Rows(174).Insert()

# Populate the newly-inserted row 174 with the new weekly price record.
$ws.Cells.Item(174, 1).Value = 2
$ws.Cells.Item(174, 2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(174, 3).Value = 'Coquimbo'
$ws.Cells.Item(174, 4).Value = 45106
$ws.Cells.Item(174, 5).Value = 4
$ws.Cells.Item(174, 6).Value = 100112024
$ws.Cells.Item(174, 7).Value = 'Choclo'
$ws.Cells.Item(174, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(174, 9).Value = 'Primera'
$ws.Cells.Item(174, 10).Value = 30000
$ws.Cells.Item(174, 11).Value = 230
$ws.Cells.Item(174, 12).Value = 250
$ws.Cells.Item(174, 13).Value = 240
$ws.Cells.Item(174, 14).Value = '$/unidad'
$ws.Cells.Item(174, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(174, 16).Value = 240
$ws.Cells.Item(174, 17).Value = 1
$ws.Cells.Item(174, 18).Value = 'Hortaliza'
